$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new client row (row 8) with the same shape as the existing rows.
$ws.Cells.Item(8, 1).Value = "Ana Milic"
$ws.Cells.Item(8, 2).Value = "067/777-888"
$ws.Cells.Item(8, 3).Value = "Tivat"

# D8/E8 are blank text cells (like D2:E7 above them). A plain "" assignment
# clears/omits the cell entirely, so use a quote-prefix entry to force an
# actual empty-string text value, then strip the quote-prefix style it adds.
$ws.Cells.Item(8, 4).Value = "'"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "'"
$ws.Cells.Item(8, 5).Style = "Normal"
